$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 0.5331505541342381
$ws.Cells.Item(2, 3).Value = 0.1057959876820505
$ws.Cells.Item(2, 5).Value = 0.09971148296677157
$ws.Cells.Item(2, 6).Value = 0.4443680307746263
$ws.Cells.Item(2, 7).Value = 0.002467003654033875
$ws.Cells.Item(2, 9).Value = 0.8606561198313081
$ws.Cells.Item(2, 11).Value = 0.3290543506840038
$ws.Cells.Item(2, 12).Value = 0.2049819237102781
$ws.Cells.Item(2, 15).Value = 3.282649908514259
$ws.Cells.Item(3, 2).Value = 0.4891433345004543
$ws.Cells.Item(3, 3).Value = 0.1049476832844576
$ws.Cells.Item(3, 5).Value = 0.09920872303841222
$ws.Cells.Item(3, 6).Value = 0.387822817061874
$ws.Cells.Item(3, 7).Value = 0.002469229783322664
$ws.Cells.Item(3, 9).Value = 0.8713952541556189
$ws.Cells.Item(3, 11).Value = 0.2908983227742965
$ws.Cells.Item(3, 12).Value = 0.1977173683576581
$ws.Cells.Item(3, 15).Value = 3.320983056847922
$ws.Cells.Item(4, 2).Value = 0.4621895096719015
$ws.Cells.Item(4, 3).Value = 0.1044315231265713
$ws.Cells.Item(4, 5).Value = 0.09895297245058998
$ws.Cells.Item(4, 6).Value = 0.3531389305168915
$ws.Cells.Item(4, 7).Value = 0.002470669215449162
$ws.Cells.Item(4, 9).Value = 0.8784481333453265
$ws.Cells.Item(4, 11).Value = 0.2674344730348821
$ws.Cells.Item(4, 12).Value = 0.1933620614794194
$ws.Cells.Item(4, 15).Value = 3.346397436567059
$ws.Cells.Item(5, 2).Value = 0.4512231724378921
$ws.Cells.Item(5, 3).Value = 0.1042223939967073
$ws.Cells.Item(5, 5).Value = 0.09886208814385355
$ws.Cells.Item(5, 6).Value = 0.3390132514313251
$ws.Cells.Item(5, 7).Value = 0.00247127409951309
$ws.Cells.Item(5, 9).Value = 0.8814376422975005
$ws.Cells.Item(5, 11).Value = 0.2578643485696404
$ws.Cells.Item(5, 12).Value = 0.1916137753094773
$ws.Cells.Item(5, 15).Value = 3.357226010167452
$ws.Cells.Item(6, 2).Value = 0.4494033041486034
$ws.Cells.Item(6, 3).Value = 0.1041877421881452
$ws.Cells.Item(6, 5).Value = 0.09884780300275153
$ws.Cells.Item(6, 6).Value = 0.336668177824194
$ws.Cells.Item(6, 7).Value = 0.002471375646755765
$ws.Cells.Item(6, 9).Value = 0.8819410175475841
$ws.Cells.Item(6, 11).Value = 0.2562747464634896
$ws.Cells.Item(6, 12).Value = 0.1913250792334225
$ws.Cells.Item(6, 15).Value = 3.359052590317347
$ws.Cells.Item(7, 2).Value = 0.4620415416674177
$ws.Cells.Item(7, 3).Value = 0.1044286977995696
$ws.Cells.Item(7, 5).Value = 0.09895169272559556
$ws.Cells.Item(7, 6).Value = 0.3529483938344953
$ws.Cells.Item(7, 7).Value = 0.002470677298900508
$ws.Cells.Item(7, 9).Value = 0.8784879836681583
$ws.Cells.Item(7, 11).Value = 0.2673054402505954
$ws.Cells.Item(7, 12).Value = 0.1933383759101304
$ws.Cells.Item(7, 15).Value = 3.346541563555093
$ws.Cells.Item(8, 2).Value = 0.517963451371827
$ws.Cells.Item(8, 3).Value = 0.1055025327852661
$ws.Cells.Item(8, 5).Value = 0.09952715543666812
$ws.Cells.Item(8, 6).Value = 0.4248636149813336
$ws.Cells.Item(8, 7).Value = 0.002467756191466366
$ws.Cells.Item(8, 9).Value = 0.8642636965744934
$ws.Cells.Item(8, 11).Value = 0.3159059561877484
$ws.Cells.Item(8, 12).Value = 0.2024553325890821
$ws.Cells.Item(8, 15).Value = 3.295477380164968
$ws.Cells.Item(9, 2).Value = 0.6281277599456701
$ws.Cells.Item(9, 3).Value = 0.1076445235300838
$ws.Cells.Item(9, 5).Value = 0.1010749362370902
$ws.Cells.Item(9, 6).Value = 0.5661985755041457
$ws.Cells.Item(9, 7).Value = 0.002462601376394921
$ws.Cells.Item(9, 9).Value = 0.8400125127313736
$ws.Cells.Item(9, 11).Value = 0.410904419114047
$ws.Cells.Item(9, 12).Value = 0.2211650685690358
$ws.Cells.Item(9, 15).Value = 3.210245566488595
$ws.Cells.Item(10, 2).Value = 0.7093414612786546
$ws.Cells.Item(10, 3).Value = 0.109239005707785
$ws.Cells.Item(10, 5).Value = 0.102466897348112
$ws.Cells.Item(10, 6).Value = 0.6702781546542269
$ws.Cells.Item(10, 7).Value = 0.002459160335301414
$ws.Cells.Item(10, 9).Value = 0.8244162953129077
$ws.Cells.Item(10, 11).Value = 0.4804896305834347
$ws.Cells.Item(10, 12).Value = 0.2354156205236819
$ws.Cells.Item(10, 15).Value = 3.156720669277647
$ws.Cells.Item(11, 2).Value = 0.7463418150489929
$ws.Cells.Item(11, 3).Value = 0.1099686215239757
$ws.Cells.Item(11, 5).Value = 0.1031553097058158
$ws.Cells.Item(11, 6).Value = 0.7176906081379002
$ws.Cells.Item(11, 7).Value = 0.002457669367557512
$ws.Cells.Item(11, 9).Value = 0.8178035771634953
$ws.Cells.Item(11, 11).Value = 0.5120956772026091
$ws.Cells.Item(11, 12).Value = 0.2420077173974562
$ws.Cells.Item(11, 15).Value = 3.134347547766382
$ws.Cells.Item(12, 2).Value = 0.7603602488313754
$ws.Cells.Item(12, 3).Value = 0.1102454978530787
$ws.Cells.Item(12, 5).Value = 0.1034239128937173
$ws.Cells.Item(12, 6).Value = 0.7356546913071611
$ws.Cells.Item(12, 7).Value = 0.002457115418087628
$ws.Cells.Item(12, 9).Value = 0.8153688595724944
$ws.Cells.Item(12, 11).Value = 0.5240565509353985
$ws.Cells.Item(12, 12).Value = 0.244519640228205
$ws.Cells.Item(12, 15).Value = 3.126159712489226
$ws.Cells.Item(13, 2).Value = 0.7573408210689934
$ws.Cells.Item(13, 3).Value = 0.1101858419414583
$ws.Cells.Item(13, 5).Value = 0.1033657126546714
$ws.Cells.Item(13, 6).Value = 0.7317853510981394
$ws.Cells.Item(13, 7).Value = 0.002457234248303131
$ws.Cells.Item(13, 9).Value = 0.8158901332440429
$ws.Cells.Item(13, 11).Value = 0.5214809134841971
$ws.Cells.Item(13, 12).Value = 0.2439779587564175
$ws.Cells.Item(13, 15).Value = 3.127910456845939
$ws.Cells.Item(14, 2).Value = 0.747494980399324
$ws.Cells.Item(14, 3).Value = 0.109991388710533
$ws.Cells.Item(14, 5).Value = 0.1031772493204137
$ws.Cells.Item(14, 6).Value = 0.7191683204515869
$ws.Cells.Item(14, 7).Value = 0.00245762358070929
$ws.Cells.Item(14, 9).Value = 0.8176018807958556
$ws.Cells.Item(14, 11).Value = 0.5130798625717432
$ws.Cells.Item(14, 12).Value = 0.2422140619490136
$ws.Cells.Item(14, 15).Value = 3.133668228578173
$ws.Cells.Item(15, 2).Value = 0.7414650322499199
$ws.Cells.Item(15, 3).Value = 0.1098723560514614
$ws.Cells.Item(15, 5).Value = 0.1030628403866807
$ws.Cells.Item(15, 6).Value = 0.7114413442032514
$ws.Cells.Item(15, 7).Value = 0.002457863443311499
$ws.Cells.Item(15, 9).Value = 0.8186594118715647
$ws.Cells.Item(15, 11).Value = 0.5079329605077589
$ws.Cells.Item(15, 12).Value = 0.2411356578874404
$ws.Cells.Item(15, 15).Value = 3.137232073262552
$ws.Cells.Item(16, 2).Value = 0.7069244539850388
$ws.Cells.Item(16, 3).Value = 0.1091914072389812
$ws.Cells.Item(16, 5).Value = 0.1024230169442255
$ws.Cells.Item(16, 6).Value = 0.6671810134426437
$ws.Cells.Item(16, 7).Value = 0.002459259266128363
$ws.Cells.Item(16, 9).Value = 0.8248581571808451
$ws.Cells.Item(16, 11).Value = 0.4784230664196798
$ws.Cells.Item(16, 12).Value = 0.2349870053608356
$ws.Cells.Item(16, 15).Value = 3.158222584611863
$ws.Cells.Item(17, 2).Value = 0.6857486397104822
$ws.Cells.Item(17, 3).Value = 0.1087747425098655
$ws.Cells.Item(17, 5).Value = 0.1020446292421013
$ws.Cells.Item(17, 6).Value = 0.6400460337125793
$ws.Cells.Item(17, 7).Value = 0.002460134574241664
$ws.Cells.Item(17, 9).Value = 0.8287843963285937
$ws.Cells.Item(17, 11).Value = 0.4603067807273646
$ws.Cells.Item(17, 12).Value = 0.2312429678412116
$ws.Cells.Item(17, 15).Value = 3.17160579305343
$ws.Cells.Item(18, 2).Value = 0.6735741759411837
$ws.Cells.Item(18, 3).Value = 0.1085354923217707
$ws.Cells.Item(18, 5).Value = 0.1018321878445931
$ws.Cells.Item(18, 6).Value = 0.6244449056556647
$ws.Cells.Item(18, 7).Value = 0.002460645032016762
$ws.Cells.Item(18, 9).Value = 0.8310880376726928
$ws.Cells.Item(18, 11).Value = 0.4498822357188885
$ws.Cells.Item(18, 12).Value = 0.2290998043675216
$ws.Cells.Item(18, 15).Value = 3.179489398344771
$ws.Cells.Item(19, 2).Value = 0.6694530479680907
$ws.Cells.Item(19, 3).Value = 0.108454556637632
$ws.Cells.Item(19, 5).Value = 0.1017611521194119
$ws.Cells.Item(19, 6).Value = 0.619163680173358
$ws.Cells.Item(19, 7).Value = 0.00246081906871754
$ws.Cells.Item(19, 9).Value = 0.8318758022006065
$ws.Cells.Item(19, 11).Value = 0.4463519048612454
$ws.Cells.Item(19, 12).Value = 0.2283759394650815
$ws.Cells.Item(19, 15).Value = 3.182190577843272
$ws.Cells.Item(20, 2).Value = 0.6880022982844309
$ws.Cells.Item(20, 3).Value = 0.1088190555409199
$ws.Cells.Item(20, 5).Value = 0.1020843715247963
$ws.Cells.Item(20, 6).Value = 0.642933953830422
$ws.Cells.Item(20, 7).Value = 0.002460040671672851
$ws.Cells.Item(20, 9).Value = 0.8283617452090191
$ws.Cells.Item(20, 11).Value = 0.4622357645759791
$ws.Cells.Item(20, 12).Value = 0.2316404609226055
$ws.Cells.Item(20, 15).Value = 3.170161881714307
$ws.Cells.Item(21, 2).Value = 0.7503867535567963
$ws.Cells.Item(21, 3).Value = 0.1100484886422777
$ws.Cells.Item(21, 5).Value = 0.1032323908952577
$ws.Cells.Item(21, 6).Value = 0.7228739723491628
$ws.Cells.Item(21, 7).Value = 0.002457508935849272
$ws.Cells.Item(21, 9).Value = 0.8170972158008496
$ws.Cells.Item(21, 11).Value = 0.5155476677352056
$ws.Cells.Item(21, 12).Value = 0.2427317375635596
$ws.Cells.Item(21, 15).Value = 3.131969311158102
$ws.Cells.Item(22, 2).Value = 0.7912002036526928
$ws.Cells.Item(22, 3).Value = 0.1108554006579183
$ws.Cells.Item(22, 5).Value = 0.1040288204446043
$ws.Cells.Item(22, 6).Value = 0.7751780083420101
$ws.Cells.Item(22, 7).Value = 0.002455916340673137
$ws.Cells.Item(22, 9).Value = 0.8101395912735825
$ws.Cells.Item(22, 11).Value = 0.550345183032789
$ws.Cells.Item(22, 12).Value = 0.2500716363759494
$ws.Cells.Item(22, 15).Value = 3.108665977573992
$ws.Cells.Item(23, 2).Value = 0.7694137627983935
$ws.Cells.Item(23, 3).Value = 0.1104244341964034
$ws.Cells.Item(23, 5).Value = 0.1035995370587166
$ws.Cells.Item(23, 6).Value = 0.7472568307830727
$ws.Cells.Item(23, 7).Value = 0.002456760677728542
$ws.Cells.Item(23, 9).Value = 0.813815986649125
$ws.Cells.Item(23, 11).Value = 0.5317774193455307
$ws.Cells.Item(23, 12).Value = 0.2461458909105119
$ws.Cells.Item(23, 15).Value = 3.120951628850179
$ws.Cells.Item(24, 2).Value = 0.6869834198234344
$ws.Cells.Item(24, 3).Value = 0.1087990206963099
$ws.Cells.Item(24, 5).Value = 0.1020663881554853
$ws.Cells.Item(24, 6).Value = 0.6416283278902171
$ws.Cells.Item(24, 7).Value = 0.002460083102573517
$ws.Cells.Item(24, 9).Value = 0.8285526812690662
$ws.Cells.Item(24, 11).Value = 0.461363699744453
$ws.Cells.Item(24, 12).Value = 0.2314607252212539
$ws.Cells.Item(24, 15).Value = 3.170814084018346
$ws.Cells.Item(25, 2).Value = 0.5982747882083288
$ws.Cells.Item(25, 3).Value = 0.1070613186822911
$ws.Cells.Item(25, 5).Value = 0.100611414942648
$ws.Cells.Item(25, 6).Value = 0.5279251897347166
$ws.Cells.Item(25, 7).Value = 0.002463934849905751
$ws.Cells.Item(25, 9).Value = 0.8461830839271052
$ws.Cells.Item(25, 11).Value = 0.3852400730747831
$ws.Cells.Item(25, 12).Value = 0.2160148219998632
$ws.Cells.Item(25, 15).Value = 3.231706480502709

